$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 was reporting on "Iene" (JPY); update it to report on "Dólar" (USD)
# with the new quotation values, matching the new report run.
$ws.Range("A2").Value = "Dólar"
$ws.Range("B2").Value = "USD"
$ws.Range("C2").Value = "$"
$ws.Range("D2").Value = 5.68
$ws.Range("E2").Value = 5.6
$ws.Range("F2").Value = "Diminuiu 1.41%"

# Report time stamp also moved forward.
$ws.Range("F5").Value = "22:08"
